# Cleaning pass over the raw "valid_values" data from the experts.
# Fixes a handful of stale / inconsistent category-value strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D37: "0.5-49 mm" clarified to "0.5mm-49 mm"
$ws.Range("D37").Value = " >1000 mm; 50 mm-999 mm; 0.5mm-49 mm; <0.4 mm"

# B3: trait label gained an " exposure" suffix
$ws.Range("B3").Value = "planktonic larval duration (PLD) exposure"

# D5: fecundity cat_vals was missing the "<1" bucket
$ws.Range("D5").Value = " <1; 1-2; 2-5; 5-10; 10-20; 20-50; 50-100; 100-1000; 1000-10,000; >10,000"

# C44/D44: category values corrected from "high; medium; low sensitivity" to "high; medium; low; none"
$ws.Range("C44").Value = "high; medium; low; none"
$ws.Range("D44").Value = "high; medium; low; none"

# C6/D6: "100+" standardized to ">100"
$ws.Range("C6").Value = "1, 2-10, 11-25, 26-50, 51-100, >100"
$ws.Range("D6").Value = "1; 2-10; 11-25; 26-50; 51-100; >100"

# Reflect the last-edited cell as the active selection, like Excel would after these edits
$ws.Range("D6").Select()
